$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 92, shifting the existing rows 92-94 down to 94-96
$ws.Range("A92:T93").EntireRow.Insert()

# Row 92 (new): Membrillo, Primera
$ws.Range("A92").Value = 8
$ws.Range("B92").Value = "Terminal La Palmera de La Serena"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 45075
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100104
$ws.Range("H92").Value = "Frutos de pepita"
$ws.Range("I92").Value = 100104003
$ws.Range("J92").Value = "Membrillo"
$ws.Range("K92").Value = "Champion"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 10
$ws.Range("N92").Value = 240000
$ws.Range("O92").Value = 250000
$ws.Range("P92").Value = 245000
$ws.Range("Q92").Value = "$/bins (450 kilos)"
$ws.Range("R92").Value = "Región de O'Higgins"
$ws.Range("S92").Value = 544
$ws.Range("T92").Value = 450

# Row 93 (new): Membrillo, Segunda
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 45075
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = "Fruta"
$ws.Range("G93").Value = 100104
$ws.Range("H93").Value = "Frutos de pepita"
$ws.Range("I93").Value = 100104003
$ws.Range("J93").Value = "Membrillo"
$ws.Range("K93").Value = "Champion"
$ws.Range("L93").Value = "Segunda"
$ws.Range("M93").Value = 16
$ws.Range("N93").Value = 210000
$ws.Range("O93").Value = 220000
$ws.Range("P93").Value = 215000
$ws.Range("Q93").Value = "$/bins (450 kilos)"
$ws.Range("R93").Value = "Región de O'Higgins"
$ws.Range("S93").Value = 478
$ws.Range("T93").Value = 450
